# Updated cryptos list on Thu Oct 10 15:51:46 UTC 2024 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns for each coin row with the
# latest scraped figures; row 15/16 additionally swap ShibaInu and
# WrappedliquidstakedEther2.0 (name/link/price/volume) to reflect their new rank order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/volume cells are stored as text (e.g. "60.624.24", "  -2.72%  "), so force
# the Text number format before writing -- otherwise Excel COM auto-coerces
# numeric-looking strings like "8.18" into real numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "60.624.24"
$ws.Range("E2").Value = "  -2.72%  "
$ws.Range("D3").Value = "2.414.99"
$ws.Range("E3").Value = "  -2.02%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "564.53"
$ws.Range("E5").Value = "  -3.08%  "
$ws.Range("D6").Value = "137.89"
$ws.Range("E6").Value = "  -3.37%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "0.534"
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "2.400.28"
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("E10").Value = "  -5.06%  "
$ws.Range("E11").Value = "  -1.19%  "
$ws.Range("D12").Value = "5.04"
$ws.Range("D14").Value = "25.71"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "0.0000166"
$ws.Range("E15").Value = "  -4.02%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "2.781.44"
$ws.Range("E16").Value = "  -4.15%  "
$ws.Range("D17").Value = "60.764.33"
$ws.Range("E17").Value = "  -2.26%  "
$ws.Range("D18").Value = "2.393.57"
$ws.Range("E18").Value = "  -2.72%  "
$ws.Range("D19").Value = "8.18"
$ws.Range("E19").Value = "  +11.27%  "
$ws.Range("D20").Value = "10.57"
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("D21").Value = "322.84"
$ws.Range("E21").Value = "  -1.36%  "
$ws.Range("D22").Value = "4.04"
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("D23").Value = "6.14"
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "1.80"
$ws.Range("E25").Value = "  -8.60%  "
$ws.Range("D26").Value = "64.23"
$ws.Range("E26").Value = "  -1.91%  "
$ws.Range("D27").Value = "552.62"
$ws.Range("E27").Value = "  -6.09%  "
$ws.Range("D28").Value = "8.05"
$ws.Range("E28").Value = "  -11.80%  "
$ws.Range("D29").Value = "2.526.57"
$ws.Range("E29").Value = "  -1.94%  "
$ws.Range("D30").Value = "0.0₃0911"
$ws.Range("E30").Value = "  -3.94%  "
$ws.Range("D31").Value = "7.85"
$ws.Range("E31").Value = "  -1.77%  "
$ws.Range("D32").Value = "1.30"
$ws.Range("E32").Value = "  -6.66%  "
$ws.Range("E33").Value = "  -5.02%  "
$ws.Range("E34").Value = "  -2.75%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").Value = "151.66"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  -1.69%  "
$ws.Range("D39").Value = "4.50"
$ws.Range("E39").Value = "  -6.32%  "
$ws.Range("D40").Value = "18.21"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("D41").Value = "5.07"
$ws.Range("E41").Value = "  -2.95%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  -4.20%  "
$ws.Range("D44").Value = "0.0₆0289"
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").Value = "2.28"
$ws.Range("E45").Value = "  -5.53%  "
$ws.Range("D46").Value = "142.70"
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("D47").Value = "3.49"
$ws.Range("E47").Value = "  -3.26%  "
$ws.Range("D48").Value = "0.583"
$ws.Range("E48").Value = "  -3.46%  "
$ws.Range("D49").Value = "0.0496"
$ws.Range("E49").Value = "  -3.63%  "
$ws.Range("E50").Value = "  -5.00%  "
$ws.Range("D51").Value = "0.0894"
$ws.Range("E51").Value = "  -0.67%  "

# Re-normalize style so no stray cell-format index lingers from the NumberFormat
# override above (matches the source file, where these cells carry no "s" attr).
$ws.Range("D2:E51").Style = "Normal"
